$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 297-298, pushing the existing rows
# (old 297..394) down to 299..396.
$ws.Rows("297:298").Insert()

# --- New row 297 ---
$ws.Cells.Item(297, 1).Value = 1
$ws.Cells.Item(297, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(297, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(297, 4).Value = 44900
$ws.Cells.Item(297, 5).Value = 15
$ws.Cells.Item(297, 6).Value = 100112043
$ws.Cells.Item(297, 7).Value = "Pepino ensalada"
$ws.Cells.Item(297, 8).Value = "Sin especificar"
$ws.Cells.Item(297, 9).Value = "Primera"
$ws.Cells.Item(297, 10).Value = 550
$ws.Cells.Item(297, 11).Value = 13000
$ws.Cells.Item(297, 12).Value = 14000
$ws.Cells.Item(297, 13).Value = 13455
$ws.Cells.Item(297, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(297, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(297, 16).Value = 192
$ws.Cells.Item(297, 17).Value = 70
$ws.Cells.Item(297, 18).Value = "Hortaliza"

# --- New row 298 ---
$ws.Cells.Item(298, 1).Value = 1
$ws.Cells.Item(298, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(298, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(298, 4).Value = 44900
$ws.Cells.Item(298, 5).Value = 15
$ws.Cells.Item(298, 6).Value = 100112043
$ws.Cells.Item(298, 7).Value = "Pepino ensalada"
$ws.Cells.Item(298, 8).Value = "Sin especificar"
$ws.Cells.Item(298, 9).Value = "Segunda"
$ws.Cells.Item(298, 10).Value = 250
$ws.Cells.Item(298, 11).Value = 10000
$ws.Cells.Item(298, 12).Value = 11000
$ws.Cells.Item(298, 13).Value = 10600
$ws.Cells.Item(298, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(298, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(298, 16).Value = 106
$ws.Cells.Item(298, 17).Value = 100
$ws.Cells.Item(298, 18).Value = "Hortaliza"
